$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Workspace Name - case sensitive"
$ws.Range("B1").Value = "Team Name(s) - comma-delimited and case sensitive (if --use_team_id parameter is passed: Team ID instead)"
$ws.Range("C1").Value = "Status (will be set to 'success' or have an error message)"

$ws.Range("C2").Select()
